# Update "Hjemme passive" data: meanEMG / legmaxROM values for columns B:E
# (corresponds to Subj 15/16 trials CON/STR) on Ark1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1: subject/header numbers
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: CON values
$ws.Range("B2").Value = 60.780713310805943
$ws.Range("C2").Value = 49.846522874245437
$ws.Range("D2").Value = 65.197813094466525
$ws.Range("E2").Value = 52.598329430698286

# Row 3: STR values
$ws.Range("B3").Value = 63.221206623705854
$ws.Range("C3").Value = 44.659973050356776
$ws.Range("D3").Value = 74.969337591465788
$ws.Range("E3").Value = 50.062774572276382

# Update the selection to match the saved view state (B1:E3)
$ws.Range("B1:E3").Select()
